# Amend corrected label annotations
# Lowercase the values in column F (the "labels" column) for all data rows,
# leaving the header row (row 1) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)  # column F
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and (-not $val.Equals(""))) {
        $lower = $val.ToLower()
        if (-not $lower.Equals($val)) {
            $cell.Value2 = $lower
        }
    }
}
